$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "nnnn/nnnnn/nnn"
$ws.Cells.Item(2, 3).Value = 0.9825396827604405
$ws.Cells.Item(2, 4).Value = 0.966666665342119
$ws.Cells.Item(2, 5).Value = 0.9962962965170542
$ws.Cells.Item(2, 6).Value = 28.567477648788028
$ws.Cells.Item(2, 7).Value = 22.5

$ws.Cells.Item(3, 2).Value = "nnnn/nnnn/nnn"
$ws.Cells.Item(3, 3).Value = 0.9787114848275812
$ws.Cells.Item(3, 4).Value = 0.9705882352941176
$ws.Cells.Item(3, 5).Value = 0.9764705896377563
$ws.Cells.Item(3, 6).Value = 28.590975259332094
$ws.Cells.Item(3, 7).Value = 21.25

$ws.Cells.Item(4, 2).Value = "nnnn/nnnnnn/nnn"
$ws.Cells.Item(4, 3).Value = 0.979831932906677
$ws.Cells.Item(4, 4).Value = 0.9627450914943919
$ws.Cells.Item(4, 5).Value = 0.9960784316062927
$ws.Cells.Item(4, 6).Value = 31.023406596277272
$ws.Cells.Item(4, 7).Value = 21.25

$ws.Cells.Item(5, 2).Value = "nnnn/nnn/nnn"
$ws.Cells.Item(5, 3).Value = 0.9752380954651605
$ws.Cells.Item(5, 4).Value = 0.9666666626930237
$ws.Cells.Item(5, 5).Value = 1.0
$ws.Cells.Item(5, 6).Value = 26.32787763039272
$ws.Cells.Item(5, 7).Value = 6.25

$ws.Cells.Item(6, 2).Value = "nnnn/nnnnnnn/nnn"
$ws.Cells.Item(6, 3).Value = 0.9873015876800295
$ws.Cells.Item(6, 4).Value = 0.9777777791023254
$ws.Cells.Item(6, 5).Value = 1.0
$ws.Cells.Item(6, 6).Value = 19.73988404675766
$ws.Cells.Item(6, 7).Value = 3.75

$ws.Cells.Item(7, 2).Value = "nnnn/nnnnnnnn/nnn"
$ws.Cells.Item(7, 3).Value = 0.9873015873015873
$ws.Cells.Item(7, 4).Value = 0.9555555383364359
$ws.Cells.Item(7, 5).Value = 1.0
$ws.Cells.Item(7, 6).Value = 30.8073609113693
$ws.Cells.Item(7, 7).Value = 3.75

$ws.Cells.Item(8, 2).Value = "nnnn/nnnnnn/nnnnnn/nnn"
$ws.Cells.Item(8, 3).Value = 0.9809523812362126
$ws.Cells.Item(8, 4).Value = 0.9666666686534882
$ws.Cells.Item(8, 5).Value = 0.9666666686534882
$ws.Cells.Item(8, 6).Value = 26.8273805161317
$ws.Cells.Item(8, 7).Value = 2.5

$ws.Cells.Item(9, 2).Value = "nnnn/nnnnn/nnnnn/nnn"
$ws.Cells.Item(9, 3).Value = 0.9857142857142858
$ws.Cells.Item(9, 4).Value = 0.9333333373069763
$ws.Cells.Item(9, 5).Value = 0.9666666686534882
$ws.Cells.Item(9, 6).Value = 35.047248200575496
$ws.Cells.Item(9, 7).Value = 2.5

$ws.Cells.Item(10, 2).Value = "nnnn/nnnnnnnnn/nnn"
$ws.Cells.Item(10, 3).Value = 0.9857142857142858
$ws.Cells.Item(10, 4).Value = 0.9833333194255829
$ws.Cells.Item(10, 5).Value = 1.0
$ws.Cells.Item(10, 6).Value = 35.1133285999298
$ws.Cells.Item(10, 7).Value = 2.5

$ws.Cells.Item(11, 2).Value = "nnnn/nnnnnnnn/nnnnnnnn/nnn"
$ws.Cells.Item(11, 3).Value = 0.9904761904761905
$ws.Cells.Item(11, 4).Value = 0.9666666388511658
$ws.Cells.Item(11, 5).Value = 0.9333333373069763
$ws.Cells.Item(11, 6).Value = 11.2289566000303
$ws.Cells.Item(11, 7).Value = 1.25

$ws.Cells.Item(12, 2).Value = "nnnn/nnnnnnnnnnnn/nnn"
$ws.Cells.Item(12, 3).Value = 0.9809523815200443
$ws.Cells.Item(12, 4).Value = 1.0
$ws.Cells.Item(12, 5).Value = 1.0
$ws.Cells.Item(12, 6).Value = 17.1249912699064
$ws.Cells.Item(12, 7).Value = 1.25

$ws.Cells.Item(13, 2).Value = "nnnn/nnnn/n/nnnn/nnnn/n/nnnn/nnn"
$ws.Cells.Item(13, 3).Value = 0.9904761910438538
$ws.Cells.Item(13, 4).Value = 1.0
$ws.Cells.Item(13, 5).Value = 1.0
$ws.Cells.Item(13, 6).Value = 17.1510696848234
$ws.Cells.Item(13, 7).Value = 1.25

$ws.Cells.Item(14, 2).Value = "nnnn/nnnnnnn/nnnnnnn/nnnnnnn/nnnnnnn/nnn"
$ws.Cells.Item(14, 3).Value = 0.9619047624724252
$ws.Cells.Item(14, 4).Value = 1.0
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 18.4992587486903
$ws.Cells.Item(14, 7).Value = 1.25

$ws.Cells.Item(15, 2).Value = "nnnn/nnnn/nnnn/nnn"
$ws.Cells.Item(15, 3).Value = 0.9904761904761905
$ws.Cells.Item(15, 4).Value = 0.9666666388511658
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 31.7194145321846
$ws.Cells.Item(15, 7).Value = 1.25

$ws.Cells.Item(16, 2).Value = "nnnn/nnnnn/n/nnnnn/nnnnn/n/nnnnn/nnn"
$ws.Cells.Item(16, 3).Value = 0.9904761904761905
$ws.Cells.Item(16, 4).Value = 0.9333333373069763
$ws.Cells.Item(16, 5).Value = 0.9333333373069763
$ws.Cells.Item(16, 6).Value = 31.7597713828087
$ws.Cells.Item(16, 7).Value = 1.25

$ws.Cells.Item(17, 2).Value = "nnnn/nnnn/nnnn/nnnnnn/nnn"
$ws.Cells.Item(17, 3).Value = 0.9809523809523809
$ws.Cells.Item(17, 4).Value = 0.9666666388511658
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 6).Value = 33.1587646325429
$ws.Cells.Item(17, 7).Value = 1.25

$ws.Cells.Item(18, 2).Value = "nnnn/nnnn/nnnn/nnnn/nnn"
$ws.Cells.Item(18, 3).Value = 0.9809523809523809
$ws.Cells.Item(18, 4).Value = 0.9333333373069763
$ws.Cells.Item(18, 5).Value = 0.9333333373069763
$ws.Cells.Item(18, 6).Value = 36.5265705664953
$ws.Cells.Item(18, 7).Value = 1.25

$ws.Cells.Item(19, 2).Value = "nnnn/nnnnnn/nnnnnn/nnnnn/nnn"
$ws.Cells.Item(19, 3).Value = 0.9809523809523809
$ws.Cells.Item(19, 4).Value = 0.9333333373069763
$ws.Cells.Item(19, 5).Value = 1.0
$ws.Cells.Item(19, 6).Value = 36.966023349762
$ws.Cells.Item(19, 7).Value = 1.25

$ws.Cells.Item(20, 2).Value = "nnnn/nnnnnnnnnn/nnn"
$ws.Cells.Item(20, 3).Value = 0.9904761904761905
$ws.Cells.Item(20, 4).Value = 0.9333333373069763
$ws.Cells.Item(20, 5).Value = 0.9333333373069763
$ws.Cells.Item(20, 6).Value = 37.7258544484774
$ws.Cells.Item(20, 7).Value = 1.25

$ws.Cells.Item(21, 2).Value = "nnnn/nnnnnnnn/nnnnnnnn/nnnnnnnn/nnnnnnnn/nnn"
$ws.Cells.Item(21, 3).Value = 0.9809523809523809
$ws.Cells.Item(21, 4).Value = 1.0
$ws.Cells.Item(21, 5).Value = 1.0
$ws.Cells.Item(21, 6).Value = 37.9731294512749
$ws.Cells.Item(21, 7).Value = 1.25

$ws.Rows(22).Delete()